# Update the multiplication problems in the table to the new set of
# generated exercises (commit "Update master to output generated at
# c8c62b6"). Each cell's "NNN×N=" text is replaced with the new
# equation via Find & Replace on $d.Content, matching whole words only
# so partial overlaps between old/new values can't cross-contaminate.

$d = $word.ActiveDocument

$d.Content.Find.Execute("847×9=", $true, $false, $false, $false, $false, $true, 1, $false, "372×7=", 2)
$d.Content.Find.Execute("116×3=", $true, $false, $false, $false, $false, $true, 1, $false, "326×2=", 2)
$d.Content.Find.Execute("320×4=", $true, $false, $false, $false, $false, $true, 1, $false, "130×9=", 2)
$d.Content.Find.Execute("706×6=", $true, $false, $false, $false, $false, $true, 1, $false, "450×8=", 2)
$d.Content.Find.Execute("994×6=", $true, $false, $false, $false, $false, $true, 1, $false, "983×8=", 2)

$d.Content.Find.Execute("184×3=", $true, $false, $false, $false, $false, $true, 1, $false, "799×9=", 2)
$d.Content.Find.Execute("990×5=", $true, $false, $false, $false, $false, $true, 1, $false, "586×3=", 2)
$d.Content.Find.Execute("833×9=", $true, $false, $false, $false, $false, $true, 1, $false, "109×5=", 2)
$d.Content.Find.Execute("838×4=", $true, $false, $false, $false, $false, $true, 1, $false, "498×7=", 2)
$d.Content.Find.Execute("289×3=", $true, $false, $false, $false, $false, $true, 1, $false, "247×9=", 2)

$d.Content.Find.Execute("535×4=", $true, $false, $false, $false, $false, $true, 1, $false, "459×8=", 2)
$d.Content.Find.Execute("112×8=", $true, $false, $false, $false, $false, $true, 1, $false, "953×5=", 2)
$d.Content.Find.Execute("529×9=", $true, $false, $false, $false, $false, $true, 1, $false, "547×3=", 2)
$d.Content.Find.Execute("998×2=", $true, $false, $false, $false, $false, $true, 1, $false, "875×5=", 2)
$d.Content.Find.Execute("723×5=", $true, $false, $false, $false, $false, $true, 1, $false, "607×8=", 2)

$d.Content.Find.Execute("609×6=", $true, $false, $false, $false, $false, $true, 1, $false, "921×8=", 2)
$d.Content.Find.Execute("671×6=", $true, $false, $false, $false, $false, $true, 1, $false, "630×9=", 2)
$d.Content.Find.Execute("401×9=", $true, $false, $false, $false, $false, $true, 1, $false, "956×7=", 2)
$d.Content.Find.Execute("647×4=", $true, $false, $false, $false, $false, $true, 1, $false, "267×2=", 2)
$d.Content.Find.Execute("914×7=", $true, $false, $false, $false, $false, $true, 1, $false, "718×6=", 2)

$d.Content.Find.Execute("701×9=", $true, $false, $false, $false, $false, $true, 1, $false, "242×8=", 2)
$d.Content.Find.Execute("612×2=", $true, $false, $false, $false, $false, $true, 1, $false, "545×7=", 2)
$d.Content.Find.Execute("732×3=", $true, $false, $false, $false, $false, $true, 1, $false, "930×7=", 2)
$d.Content.Find.Execute("417×3=", $true, $false, $false, $false, $false, $true, 1, $false, "124×7=", 2)
$d.Content.Find.Execute("878×9=", $true, $false, $false, $false, $false, $true, 1, $false, "332×2=", 2)
